# Replace the worked long-division answers in each table cell with the
# newly generated problems/answers, one Find/Replace per cell.
# Find.Execute signature used:
#   (FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
#    MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)
# Wrap = 1 (wdFindContinue), Replace = 2 (wdReplaceAll)
$d = $word.ActiveDocument

$d.Content.Find.Execute('14÷9=1, 5', $true, $false, $false, $false, $false, $true, 1, $false, '30÷8=3, 6', 2) | Out-Null
$d.Content.Find.Execute('37÷9=4, 1', $true, $false, $false, $false, $false, $true, 1, $false, '11÷7=1, 4', 2) | Out-Null
$d.Content.Find.Execute('80÷5=16, 0', $true, $false, $false, $false, $false, $true, 1, $false, '80÷3=26, 2', 2) | Out-Null
$d.Content.Find.Execute('91÷9=10, 1', $true, $false, $false, $false, $false, $true, 1, $false, '45÷3=15, 0', 2) | Out-Null
$d.Content.Find.Execute('53÷9=5, 8', $true, $false, $false, $false, $false, $true, 1, $false, '46÷5=9, 1', 2) | Out-Null
$d.Content.Find.Execute('64÷4=16, 0', $true, $false, $false, $false, $false, $true, 1, $false, '18÷7=2, 4', 2) | Out-Null
$d.Content.Find.Execute('44÷3=14, 2', $true, $false, $false, $false, $false, $true, 1, $false, '56÷3=18, 2', 2) | Out-Null
$d.Content.Find.Execute('63÷9=7, 0', $true, $false, $false, $false, $false, $true, 1, $false, '20÷2=10, 0', 2) | Out-Null
$d.Content.Find.Execute('40÷4=10, 0', $true, $false, $false, $false, $false, $true, 1, $false, '83÷5=16, 3', 2) | Out-Null
$d.Content.Find.Execute('79÷8=9, 7', $true, $false, $false, $false, $false, $true, 1, $false, '50÷2=25, 0', 2) | Out-Null
$d.Content.Find.Execute('91÷7=13, 0', $true, $false, $false, $false, $false, $true, 1, $false, '92÷8=11, 4', 2) | Out-Null
$d.Content.Find.Execute('70÷8=8, 6', $true, $false, $false, $false, $false, $true, 1, $false, '88÷6=14, 4', 2) | Out-Null
$d.Content.Find.Execute('84÷2=42, 0', $true, $false, $false, $false, $false, $true, 1, $false, '91÷6=15, 1', 2) | Out-Null
$d.Content.Find.Execute('90÷9=10, 0', $true, $false, $false, $false, $false, $true, 1, $false, '19÷7=2, 5', 2) | Out-Null
$d.Content.Find.Execute('27÷9=3, 0', $true, $false, $false, $false, $false, $true, 1, $false, '43÷7=6, 1', 2) | Out-Null
$d.Content.Find.Execute('39÷6=6, 3', $true, $false, $false, $false, $false, $true, 1, $false, '32÷9=3, 5', 2) | Out-Null
$d.Content.Find.Execute('65÷2=32, 1', $true, $false, $false, $false, $false, $true, 1, $false, '78÷7=11, 1', 2) | Out-Null
$d.Content.Find.Execute('39÷7=5, 4', $true, $false, $false, $false, $false, $true, 1, $false, '98÷2=49, 0', 2) | Out-Null
$d.Content.Find.Execute('33÷9=3, 6', $true, $false, $false, $false, $false, $true, 1, $false, '24÷4=6, 0', 2) | Out-Null
$d.Content.Find.Execute('94÷7=13, 3', $true, $false, $false, $false, $false, $true, 1, $false, '89÷8=11, 1', 2) | Out-Null
$d.Content.Find.Execute('12÷8=1, 4', $true, $false, $false, $false, $false, $true, 1, $false, '37÷8=4, 5', 2) | Out-Null
$d.Content.Find.Execute('47÷5=9, 2', $true, $false, $false, $false, $false, $true, 1, $false, '38÷4=9, 2', 2) | Out-Null
$d.Content.Find.Execute('42÷8=5, 2', $true, $false, $false, $false, $false, $true, 1, $false, '42÷6=7, 0', 2) | Out-Null
$d.Content.Find.Execute('42÷5=8, 2', $true, $false, $false, $false, $false, $true, 1, $false, '67÷6=11, 1', 2) | Out-Null
$d.Content.Find.Execute('35÷8=4, 3', $true, $false, $false, $false, $false, $true, 1, $false, '72÷6=12, 0', 2) | Out-Null
